$d = $word.ActiveDocument

# --- Edit 1: Update the first paragraph ---
# "This is a Microsoft word document." becomes:
#   "This is a Microsoft word document.  " (plain run, 2 trailing spaces)
#   + 3 red runs: "(This is a change – Ve" | "rsion for main branch" | ")"
$p1 = $d.Paragraphs(1)
$rng = $p1.Range
$rng.End = $rng.End - 1
$insertStart = $rng.End
$rng.InsertAfter("  (This is a change – Version for main branch)")
$insertEnd = $rng.End

$redStart = $insertStart + 2

$r1 = $d.Range($redStart, $redStart + 22)
$r1.Font.Color = 255

$r2 = $d.Range($redStart + 22, $redStart + 43)
$r2.Font.Color = 255

$r3 = $d.Range($redStart + 43, $insertEnd)
$r3.Font.Color = 255

# --- Edit 2: Append a new, empty, shaded paragraph at the end of the document ---
$endRng = $d.Range($d.Content.End, $d.Content.End)
$xml = @'
<?xml version="1.0" standalone="yes"?>
<?mso-application progid="Word.Document"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/_rels/.rels" pkg:contentType="application/vnd.openxmlformats-package.relationships+xml" pkg:padding="512"><pkg:xmlData><Relationships xmlns="http://schemas.openxmlformats.org/package/2006/relationships"><Relationship Id="rId1" Type="http://schemas.openxmlformats.org/officeDocument/2006/relationships/officeDocument" Target="word/document.xml"/></Relationships></pkg:xmlData></pkg:part><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:shd w:val="clear" w:color="auto" w:fill="F9F9F9"/></w:pPr></w:p><w:sectPr><w:pgSz w:w="12240" w:h="15840"/><w:pgMar w:top="1440" w:right="1440" w:bottom="1440" w:left="1440" w:header="720" w:footer="720" w:gutter="0"/><w:cols w:space="720"/></w:sectPr></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$endRng.InsertXML($xml)
